$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Anthoni'
$ws.Range("B2").Value = 'Goleiro'
$ws.Range("C2").Value = 'INT'
$ws.Range("D2").Value = 32

$ws.Range("A3").Value = 'Ronaldo'
$ws.Range("B3").Value = 'Goleiro'
$ws.Range("C3").Value = 'BAH'
$ws.Range("D3").Value = 31

$ws.Range("A4").Value = 'Weverton'
$ws.Range("B4").Value = 'Goleiro'
$ws.Range("C4").Value = 'PAL'
$ws.Range("D4").Value = 27

$ws.Range("A5").Value = 'Léo Ortiz'
$ws.Range("B5").Value = 'Zagueiro'
$ws.Range("C5").Value = 'FLA'
$ws.Range("D5").Value = 82

$ws.Range("A6").Value = 'Junior Alonso'
$ws.Range("B6").Value = 'Zagueiro'
$ws.Range("C6").Value = 'CAM'
$ws.Range("D6").Value = 29

$ws.Range("A7").Value = 'Jair'
$ws.Range("B7").Value = 'Zagueiro'
$ws.Range("C7").Value = 'BOT'
$ws.Range("D7").Value = 24

$ws.Range("A8").Value = 'Vitão'
$ws.Range("B8").Value = 'Zagueiro'
$ws.Range("C8").Value = 'INT'
$ws.Range("D8").Value = 20

$ws.Range("A9").Value = 'Bernabéi'
$ws.Range("B9").Value = 'Lateral'
$ws.Range("C9").Value = 'INT'
$ws.Range("D9").Value = 55

$ws.Range("A10").Value = 'Juninho Capixaba'
$ws.Range("B10").Value = 'Lateral'
$ws.Range("C10").Value = 'RBB'
$ws.Range("D10").Value = 53

$ws.Range("A11").Value = 'Wesley'
$ws.Range("B11").Value = 'Lateral'
$ws.Range("C11").Value = 'FLA'
$ws.Range("D11").Value = 27

$ws.Range("A12").Value = 'Angileri'
$ws.Range("B12").Value = 'Lateral'
$ws.Range("C12").Value = 'COR'
$ws.Range("D12").Value = 27

$ws.Range("A13").Value = 'Arias'
$ws.Range("B13").Value = 'Meia'
$ws.Range("C13").Value = 'FLU'
$ws.Range("D13").Value = 74

$ws.Range("A14").Value = 'Arrascaeta'
$ws.Range("B14").Value = 'Meia'
$ws.Range("C14").Value = 'FLA'
$ws.Range("D14").Value = 74

$ws.Range("A15").Value = 'Alan Patrick'
$ws.Range("B15").Value = 'Meia'
$ws.Range("C15").Value = 'INT'
$ws.Range("D15").Value = 68

$ws.Range("A16").Value = 'Gustavo Scarpa'
$ws.Range("B16").Value = 'Meia'
$ws.Range("C16").Value = 'CAM'
$ws.Range("D16").Value = 29

$ws.Range("A17").Value = 'Gerson'
$ws.Range("B17").Value = 'Meia'
$ws.Range("C17").Value = 'FLA'
$ws.Range("D17").Value = 24

$ws.Range("A18").Value = 'Yuri Alberto'
$ws.Range("B18").Value = 'Atacante'
$ws.Range("C18").Value = 'COR'
$ws.Range("D18").Value = 76

$ws.Range("A19").Value = 'Guilherme'
$ws.Range("B19").Value = 'Atacante'
$ws.Range("C19").Value = 'SAN'
$ws.Range("D19").Value = 47

$ws.Range("A20").Value = 'Igor Jesus'
$ws.Range("B20").Value = 'Atacante'
$ws.Range("C20").Value = 'BOT'
$ws.Range("D20").Value = 41

$ws.Range("A21").Value = 'Hulk'
$ws.Range("B21").Value = 'Atacante'
$ws.Range("C21").Value = 'CAM'
$ws.Range("D21").Value = 36

$ws.Range("A22").Value = 'Estêvão'
$ws.Range("B22").Value = 'Atacante'
$ws.Range("C22").Value = 'PAL'
$ws.Range("D22").Value = 35
